# Update DB connections and arrange all code
# Applies updated branch-wise stock status figures to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 30
$ws.Range("F2").Value = 1

# Row 3
$ws.Range("D3").Value = 34
$ws.Range("E3").Value = 12
$ws.Range("F3").Value = 8
$ws.Range("G3").Value = 13
$ws.Range("H3").Value = 39

# Row 4
$ws.Range("D4").Value = 38
$ws.Range("G4").Value = 8

# Row 5
$ws.Range("D5").Value = 41
$ws.Range("E5").Value = 12

# Row 6
$ws.Range("D6").Value = 39
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 12

# Row 7
$ws.Range("G7").Value = 18
$ws.Range("H7").Value = 35

# Row 8
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 17
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 12
$ws.Range("H8").Value = 44

# Row 9
$ws.Range("F9").Value = 12
$ws.Range("G9").Value = 5

# Row 10
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 4

# Row 11
$ws.Range("D11").Value = 36
$ws.Range("E11").Value = 16
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 12

# Row 12
$ws.Range("D12").Value = 38
$ws.Range("E12").Value = 7

# Row 13
$ws.Range("E13").Value = 13
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 12
$ws.Range("H13").Value = 37

# Row 14
$ws.Range("E14").Value = 21
$ws.Range("F14").Value = 7
$ws.Range("H14").Value = 30

# Row 15
$ws.Range("G15").Value = 8
$ws.Range("H15").Value = 42

# Row 16
$ws.Range("D16").Value = 32
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 52

# Row 17
$ws.Range("E17").Value = 13
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 9

# Row 18
$ws.Range("D18").Value = 44
$ws.Range("E18").Value = 19
$ws.Range("H18").Value = 30

# Row 19
$ws.Range("D19").Value = 33
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 5
$ws.Range("H19").Value = 55

# Row 20
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 3
$ws.Range("H20").Value = 57

# Row 21
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = 15
$ws.Range("F21").Value = 11
$ws.Range("G21").Value = 15
$ws.Range("H21").Value = 29

# Row 22
$ws.Range("D22").Value = 36
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 5

# Row 24
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = 20
$ws.Range("F24").Value = 6
$ws.Range("G24").Value = 11
$ws.Range("H24").Value = 26

# Row 25
$ws.Range("E25").Value = 12
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 38

# Row 27
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 42

# Row 28
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 14
$ws.Range("F28").Value = 7

# Row 29
$ws.Range("D29").Value = 37
$ws.Range("E29").Value = 11
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 9
$ws.Range("H29").Value = 43

# Row 30
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 13
$ws.Range("H30").Value = 48

# Row 31
$ws.Range("E31").Value = 17
$ws.Range("F31").Value = 12

# Row 32
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 6
